$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 93.34964983437608
$ws.Cells.Item(2, 5).Value = 96.34729766845705
$ws.Cells.Item(2, 6).Value = 99.6455315083836
$ws.Cells.Item(2, 7).Value = 89.36090117946195
$ws.Cells.Item(2, 8).Value = 112434397
$ws.Cells.Item(2, 9).Value = "ROK"

$ws.Cells.Item(3, 4).Value = 102.8917147071795
$ws.Cells.Item(3, 5).Value = 95.37022399902344
$ws.Cells.Item(3, 6).Value = 103.5287129866631
$ws.Cells.Item(3, 7).Value = 93.63072821807478
$ws.Cells.Item(3, 8).Value = 112434397
$ws.Cells.Item(3, 9).Value = "ROK"

$ws.Cells.Item(4, 4).Value = 83.29361067611514
$ws.Cells.Item(4, 5).Value = 89.65027618408203
$ws.Cells.Item(4, 6).Value = 90.72614315905876
$ws.Cells.Item(4, 7).Value = 80.87085479227751
$ws.Cells.Item(4, 8).Value = 112434397
$ws.Cells.Item(4, 9).Value = "ROK"

$ws.Cells.Item(5, 4).Value = 83.27975921516442
$ws.Cells.Item(5, 5).Value = 79.02935791015625
$ws.Cells.Item(5, 6).Value = 84.29688149031561
$ws.Cells.Item(5, 7).Value = 72.38086871919515
$ws.Cells.Item(5, 8).Value = 112434397
$ws.Cells.Item(5, 9).Value = "ROK"

$ws.Cells.Item(6, 4).Value = 93.53217962863602
$ws.Cells.Item(6, 5).Value = 94.54874420166016
$ws.Cells.Item(6, 6).Value = 98.48167963170189
$ws.Cells.Item(6, 7).Value = 91.18241822098348
$ws.Cells.Item(6, 8).Value = 112434397
$ws.Cells.Item(6, 9).Value = "ROK"

$ws.Cells.Item(7, 4).Value = 96.04321793788854
$ws.Cells.Item(7, 5).Value = 95.9342041015625
$ws.Cells.Item(7, 6).Value = 101.2592219471535
$ws.Cells.Item(7, 7).Value = 92.99076654174848
$ws.Cells.Item(7, 9).Value = "ROK"

$ws.Cells.Item(8, 4).Value = 102.6037774558616
$ws.Cells.Item(8, 5).Value = 101.009162902832
$ws.Cells.Item(8, 6).Value = 106.4342292381842
$ws.Cells.Item(8, 7).Value = 96.57123768373448
$ws.Cells.Item(8, 8).Value = 112434397
$ws.Cells.Item(8, 9).Value = "ROK"

$ws.Cells.Item(9, 4).Value = 115.3616583268855
$ws.Cells.Item(9, 5).Value = 125.6245269775391
$ws.Cells.Item(9, 6).Value = 131.9486185233716
$ws.Cells.Item(9, 7).Value = 114.7504696694308
$ws.Cells.Item(9, 8).Value = 112434397
$ws.Cells.Item(9, 9).Value = "ROK"

$ws.Cells.Item(10, 4).Value = 132.592493303619
$ws.Cells.Item(10, 5).Value = 134.2390289306641
$ws.Cells.Item(10, 6).Value = 140.7654204918692
$ws.Cells.Item(10, 7).Value = 126.5267828453548
$ws.Cells.Item(10, 8).Value = 112434397
$ws.Cells.Item(10, 9).Value = "ROK"

$ws.Cells.Item(11, 4).Value = 139.4584195667709
$ws.Cells.Item(11, 5).Value = 141.4729766845703
$ws.Cells.Item(11, 6).Value = 143.9933117746886
$ws.Cells.Item(11, 7).Value = 135.4464671968856
$ws.Cells.Item(11, 8).Value = 112434397
$ws.Cells.Item(11, 9).Value = "ROK"

$ws.Cells.Item(12, 4).Value = 153.9839547775121
$ws.Cells.Item(12, 5).Value = 172.94775390625
$ws.Cells.Item(12, 6).Value = 181.4737057319477
$ws.Cells.Item(12, 7).Value = 153.2949859204586
$ws.Cells.Item(12, 8).Value = 112434397
$ws.Cells.Item(12, 9).Value = "ROK"

$ws.Cells.Item(13, 4).Value = 171.1082578223164
$ws.Cells.Item(13, 5).Value = 170.658447265625
$ws.Cells.Item(13, 6).Value = 181.11646681414
$ws.Cells.Item(13, 7).Value = 168.8678768172493
$ws.Cells.Item(13, 8).Value = 112434397
$ws.Cells.Item(13, 9).Value = "ROK"

$ws.Cells.Item(14, 4).Value = 151.5154845359957
$ws.Cells.Item(14, 5).Value = 142.9570007324219
$ws.Cells.Item(14, 6).Value = 155.4254493331076
$ws.Cells.Item(14, 7).Value = 135.3803567760394
$ws.Cells.Item(14, 8).Value = 112434397
$ws.Cells.Item(14, 9).Value = "ROK"

$ws.Cells.Item(15, 4).Value = 143.4463601924426
$ws.Cells.Item(15, 5).Value = 163.8139190673828
$ws.Cells.Item(15, 6).Value = 164.5126375408322
$ws.Cells.Item(15, 7).Value = 143.2367406523238
$ws.Cells.Item(15, 8).Value = 112434397
$ws.Cells.Item(15, 9).Value = "ROK"

$ws.Cells.Item(16, 4).Value = 165.6848318867768
$ws.Cells.Item(16, 5).Value = 144.6307067871094
$ws.Cells.Item(16, 6).Value = 168.5382879011602
$ws.Cells.Item(16, 7).Value = 134.9201786170676
$ws.Cells.Item(16, 8).Value = 112434397
$ws.Cells.Item(16, 9).Value = "ROK"

$ws.Cells.Item(17, 4).Value = 131.1091752019178
$ws.Cells.Item(17, 5).Value = 149.6473693847656
$ws.Cells.Item(17, 6).Value = 158.8105247783891
$ws.Cells.Item(17, 7).Value = 128.1872048454433
$ws.Cells.Item(17, 8).Value = 112434397
$ws.Cells.Item(17, 9).Value = "ROK"

$ws.Cells.Item(18, 4).Value = 157.4783235959013
$ws.Cells.Item(18, 5).Value = 160.4075775146484
$ws.Cells.Item(18, 6).Value = 169.9675785293012
$ws.Cells.Item(18, 7).Value = 155.605377583582
$ws.Cells.Item(18, 8).Value = 112434397
$ws.Cells.Item(18, 9).Value = "ROK"

$ws.Cells.Item(19, 4).Value = 148.2454765404049
$ws.Cells.Item(19, 5).Value = 143.5405426025391
$ws.Cells.Item(19, 6).Value = 151.5666047762299
$ws.Cells.Item(19, 7).Value = 136.041228689653
$ws.Cells.Item(19, 8).Value = 112434397
$ws.Cells.Item(19, 9).Value = "ROK"

$ws.Cells.Item(20, 4).Value = 149.4581417332705
$ws.Cells.Item(20, 5).Value = 154.5068664550781
$ws.Cells.Item(20, 6).Value = 161.091758471209
$ws.Cells.Item(20, 7).Value = 137.0968790680699
$ws.Cells.Item(20, 8).Value = 112434397
$ws.Cells.Item(20, 9).Value = "ROK"

$ws.Cells.Item(21, 4).Value = 184.1562900396977
$ws.Cells.Item(21, 5).Value = 173.16943359375
$ws.Cells.Item(21, 6).Value = 187.8788049474516
$ws.Cells.Item(21, 7).Value = 171.4527295190611
$ws.Cells.Item(21, 8).Value = 112434397
$ws.Cells.Item(21, 9).Value = "ROK"

$ws.Cells.Item(22, 4).Value = 129.8261293281296
$ws.Cells.Item(22, 5).Value = 172.0722808837891
$ws.Cells.Item(22, 6).Value = 182.1434421329816
$ws.Cells.Item(22, 7).Value = 127.1925468857417
$ws.Cells.Item(22, 8).Value = 112434397
$ws.Cells.Item(22, 9).Value = "ROK"

$ws.Cells.Item(23, 4).Value = 194.141256523824
$ws.Cells.Item(23, 5).Value = 199.1627807617188
$ws.Cells.Item(23, 6).Value = 210.8401009960508
$ws.Cells.Item(23, 7).Value = 189.9414306795937
$ws.Cells.Item(23, 8).Value = 112434397
$ws.Cells.Item(23, 9).Value = "ROK"

$ws.Cells.Item(24, 4).Value = 204.1887661460069
$ws.Cells.Item(24, 5).Value = 217.4394378662109
$ws.Cells.Item(24, 6).Value = 231.038574708744
$ws.Cells.Item(24, 7).Value = 197.6230306388024
$ws.Cells.Item(24, 8).Value = 112434397
$ws.Cells.Item(24, 9).Value = "ROK"

$ws.Cells.Item(25, 4).Value = 232.0397492809109
$ws.Cells.Item(25, 5).Value = 228.9445343017578
$ws.Cells.Item(25, 6).Value = 247.7184881499069
$ws.Cells.Item(25, 7).Value = 219.6957450265134
$ws.Cells.Item(25, 8).Value = 112434397
$ws.Cells.Item(25, 9).Value = "ROK"

$ws.Cells.Item(26, 4).Value = 246.6967286161554
$ws.Cells.Item(26, 5).Value = 244.4947357177734
$ws.Cells.Item(26, 6).Value = 252.1646665132017
$ws.Cells.Item(26, 7).Value = 237.9350342237412
$ws.Cells.Item(26, 8).Value = 112434397
$ws.Cells.Item(26, 9).Value = "ROK"

$ws.Cells.Item(27, 4).Value = 267.5456962150533
$ws.Cells.Item(27, 5).Value = 285.5765380859375
$ws.Cells.Item(27, 6).Value = 286.8491893307238
$ws.Cells.Item(27, 7).Value = 260.4299677913896
$ws.Cells.Item(27, 8).Value = 112434397
$ws.Cells.Item(27, 9).Value = "ROK"

$ws.Cells.Item(28, 4).Value = 273.6709532146991
$ws.Cells.Item(28, 5).Value = 297.7095642089844
$ws.Cells.Item(28, 6).Value = 298.8746768025113
$ws.Cells.Item(28, 7).Value = 269.8959997897244
$ws.Cells.Item(28, 8).Value = 112434397
$ws.Cells.Item(28, 9).Value = "ROK"

$ws.Cells.Item(29, 4).Value = 326.2681218351806
$ws.Cells.Item(29, 5).Value = 270.4981079101562
$ws.Cells.Item(29, 6).Value = 326.4084062727224
$ws.Cells.Item(29, 7).Value = 257.5446453911992
$ws.Cells.Item(29, 8).Value = 112434397
$ws.Cells.Item(29, 9).Value = "ROK"

$ws.Cells.Item(30, 4).Value = 264.8712834875788
$ws.Cells.Item(30, 5).Value = 237.2723083496093
$ws.Cells.Item(30, 6).Value = 267.3973573917381
$ws.Cells.Item(30, 7).Value = 236.5022865361657
$ws.Cells.Item(30, 8).Value = 112434397
$ws.Cells.Item(30, 9).Value = "ROK"

$ws.Cells.Item(31, 4).Value = 188.4909080153895
$ws.Cells.Item(31, 5).Value = 241.0477905273437
$ws.Cells.Item(31, 6).Value = 241.6237827713452
$ws.Cells.Item(31, 7).Value = 182.1455560319309
$ws.Cells.Item(31, 8).Value = 112434397
$ws.Cells.Item(31, 9).Value = "ROK"

$ws.Cells.Item(32, 4).Value = 208.4215674045494
$ws.Cells.Item(32, 5).Value = 242.138916015625
$ws.Cells.Item(32, 6).Value = 246.4733183332733
$ws.Cells.Item(32, 7).Value = 206.3729150775578
$ws.Cells.Item(32, 8).Value = 112434397
$ws.Cells.Item(32, 9).Value = "ROK"

$ws.Cells.Item(33, 4).Value = 247.5118521128308
$ws.Cells.Item(33, 5).Value = 268.7421264648437
$ws.Cells.Item(33, 6).Value = 280.3482575366564
$ws.Cells.Item(33, 7).Value = 243.4620912986651
$ws.Cells.Item(33, 8).Value = 112434397
$ws.Cells.Item(33, 9).Value = "ROK"

$ws.Cells.Item(34, 4).Value = 278.5624625449708
$ws.Cells.Item(34, 5).Value = 271.1384582519531
$ws.Cells.Item(34, 6).Value = 281.3655876518319
$ws.Cells.Item(34, 7).Value = 258.1273272421927
$ws.Cells.Item(34, 8).Value = 112434397
$ws.Cells.Item(34, 9).Value = "ROK"

$ws.Cells.Item(35, 4).Value = 316.5659588766714
$ws.Cells.Item(35, 5).Value = 323.1384582519531
$ws.Cells.Item(35, 6).Value = 334.890151530689
$ws.Cells.Item(35, 7).Value = 309.0325379412954
$ws.Cells.Item(35, 8).Value = 112434397
$ws.Cells.Item(35, 9).Value = "ROK"

$ws.Cells.Item(36, 4).Value = 275.6136757108952
$ws.Cells.Item(36, 5).Value = 253.5317687988281
$ws.Cells.Item(36, 6).Value = 292.4186645647633
$ws.Cells.Item(36, 7).Value = 251.4094256938886
$ws.Cells.Item(36, 8).Value = 112434397
$ws.Cells.Item(36, 9).Value = "ROK"

$ws.Cells.Item(37, 4).Value = 297.8800282546317
$ws.Cells.Item(37, 5).Value = 245.5319366455078
$ws.Cells.Item(37, 6).Value = 300.5168222005261
$ws.Cells.Item(37, 7).Value = 244.3977297690179
$ws.Cells.Item(37, 8).Value = 112434397
$ws.Cells.Item(37, 9).Value = "ROK"

$ws.Cells.Item(38, 4).Value = 283.4833811503157
$ws.Cells.Item(38, 5).Value = 263.8431396484375
$ws.Cells.Item(38, 6).Value = 283.5418028536292
$ws.Cells.Item(38, 7).Value = 262.3922844894456
$ws.Cells.Item(38, 8).Value = 112434397
$ws.Cells.Item(38, 9).Value = "ROK"

$ws.Cells.Item(39, 4).Value = 268.844145926903
$ws.Cells.Item(39, 5).Value = 272.5908508300781
$ws.Cells.Item(39, 6).Value = 288.2429366323893
$ws.Cells.Item(39, 7).Value = 252.4094379220857
$ws.Cells.Item(39, 8).Value = 112434397
$ws.Cells.Item(39, 9).Value = "ROK"

$ws.Cells.Item(40, 4).Value = 263.6840048298125
$ws.Cells.Item(40, 5).Value = 262.1799926757812
$ws.Cells.Item(40, 6).Value = 272.3443669061796
$ws.Cells.Item(40, 7).Value = 255.8887009856912
$ws.Cells.Item(40, 8).Value = 112434397
$ws.Cells.Item(40, 9).Value = "ROK"

$ws.Cells.Item(41, 4).Value = 283.3833999859353
$ws.Cells.Item(41, 5).Value = 274.959716796875
$ws.Cells.Item(41, 6).Value = 289.3283876272331
$ws.Cells.Item(41, 7).Value = 266.1805350159948
$ws.Cells.Item(41, 8).Value = 112434397
$ws.Cells.Item(41, 9).Value = "ROK"

$ws.Cells.Item(42, 4).Value = 255.5065191869746
$ws.Cells.Item(42, 5).Value = 245.6862182617188
$ws.Cells.Item(42, 6).Value = 260.505956670928
$ws.Cells.Item(42, 7).Value = 213.2692929921693
$ws.Cells.Item(42, 8).Value = 112434397
$ws.Cells.Item(42, 9).Value = "ROK"

$ws.Cells.Item(43, 4).Value = 329.1354538089378
$ws.Cells.Item(43, 5).Value = 350.3638916015625
$ws.Cells.Item(43, 6).Value = 359.5386640518012
$ws.Cells.Item(43, 7).Value = 328.6672514241782
$ws.Cells.Item(43, 8).Value = 112434397
$ws.Cells.Item(43, 9).Value = "ROK"
